$d = $word.ActiveDocument

# Remove the existing hidden "_GoBack" bookmark so we can re-anchor it later.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Insert the two new paragraphs after the first paragraph's text, with the
# "_GoBack" bookmark re-created at the end of the last new paragraph.
$insertAt = $d.Paragraphs(1).Range.End - 1
$r = $d.Range($insertAt, $insertAt)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Making some changes</w:t></w:r></w:p><w:p><w:r><w:t>Hopefully this change will stay</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)
